# "Add Member From Admin Panel Fixed"
# Updates the Users report:
#  - Row 3 (Mohamed G) gets a fresh last_login / Events Joined / Participation %
#  - Row 7 (previously the "Tyra" test row) is replaced by a new "New Form" pending signup
#  - A new Row 8 is appended for member "Alan 2"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the literal string into the cell without Excel re-interpreting
    # look-alike dates/percentages as numbers, then drop back to the
    # workbook's default (unstyled) cell format.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Row 3: Mohamed G -- refresh last_login / participation -----------------
$ws.Cells.Item(3, 9).Value = "2024-12-09 01:54:45"   # last_login
$ws.Cells.Item(3, 10).Value = 1                        # Events Joined
Set-TextValue $ws.Cells.Item(3, 11) "20.0%"            # Participation Percentage

# --- Row 7: replace the old "Tyra" row with the new "New Form" signup -------
$ws.Cells.Item(7, 1).Value = 14
$ws.Cells.Item(7, 2).Value = "New Form"
$ws.Cells.Item(7, 3).Value = "12345@gmail.com"
$ws.Cells.Item(7, 4).Value = "pending"
$ws.Cells.Item(7, 5).Value = "none"
Set-TextValue $ws.Cells.Item(7, 6) "09/02/2000"
$ws.Cells.Item(7, 7).Value = "Cambridge"
$ws.Cells.Item(7, 8).Value = "2024-12-09 02:43:59"
$ws.Cells.Item(7, 9).Value = "2024-12-09 02:43:59"
$ws.Cells.Item(7, 10).Value = 3
Set-TextValue $ws.Cells.Item(7, 11) "60.0%"

# --- Row 8: new member "Alan 2" ---------------------------------------------
$ws.Cells.Item(8, 1).Value = 17
$ws.Cells.Item(8, 2).Value = "Alan 2"
$ws.Cells.Item(8, 3).Value = "alan2@works.com"
$ws.Cells.Item(8, 4).Value = "active"
$ws.Cells.Item(8, 5).Value = "none"
# Date of Birth / City are blank for this member, but keep the cells present
# (unstyled) so row 8 has the same shape as the other rows.
$ws.Cells.Item(8, 6).NumberFormat = "General"
$ws.Cells.Item(8, 6).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "General"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(8, 8).Value = "2024-12-09 03:07:23"
$ws.Cells.Item(8, 9).Value = "2024-12-09 03:07:23"
$ws.Cells.Item(8, 10).Value = 0
Set-TextValue $ws.Cells.Item(8, 11) "0.0%"
